$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("begroting")

# --- B13 gets style applied (same style as e.g. B7, "Good" named style) ---
$ws.Range("B13").Style = "Good"

# --- Restructure rows 27-38 into rows 26-42 ---
# First, clear old range B27:C38 area that will be rebuilt (A27:C38)
$ws.Range("A26:I42").ClearContents()

# Row 26: label "hours spent"
$ws.Range("A26").Value = "hours spent"

# Rows 27-37: index numbers in col A, hours in col B
$ws.Range("A27").Value = 2
$ws.Range("B27").Value = 36

$ws.Range("A28").Value = 3
$ws.Range("B28").Value = 36

$ws.Range("A29").Value = 4
$ws.Range("B29").Value = 20
$ws.Range("C29").NumberFormat = "0.0"

$ws.Range("A30").Value = 5
$ws.Range("B30").Value = 31

$ws.Range("A31").Value = 6
$ws.Range("B31").Value = 20

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = 26

$ws.Range("A33").Value = 8
$ws.Range("B33").Value = 20

$ws.Range("A34").Value = 9
$ws.Range("B34").Value = 22

$ws.Range("A35").Value = 10
$ws.Range("B35").Value = 14

$ws.Range("A36").Value = 13
$ws.Range("B36").Value = 30

$ws.Range("A37").Value = 14
$ws.Range("B37").Formula = "=14+19"

# Row 42: totals
$ws.Range("A42").Value = "Total hours spent so far"
$ws.Range("B42").Formula = "=SUM(B27:B40)"
$ws.Range("C42").Formula = "=B42*135"

# --- Sheet view: topLeftCell = A10 ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1

# --- Window sizing (workbookView) ---
$excel.Width = 24270
$excel.Height = 15990
$excel.Left = 1050
$excel.Top = -120
